$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataElement")

# Fill in the HL7 v2 mapping values (columns B, C, D) for rows that were
# previously blank. Each row corresponds to a DataElement.* field listed
# in column A. Column B/C are filled in row order first, then column D,
# matching the order the values were actually typed in.

$ws.Range("B4").Value = "OM1.2"
$ws.Range("C4").Value = "OM1.2"

$ws.Range("B6").Value = "MFE.1"
$ws.Range("C6").Value = "MFE.1"

$ws.Range("B8").Value = "OM1.21"
$ws.Range("C8").Value = "OM1.21"

$ws.Range("B9").Value = "OM1.16"
$ws.Range("C9").Value = "OM1.16"

$ws.Range("B11").Value = "OM1.11"
$ws.Range("C11").Value = "OM1.11"

$ws.Range("B12").Value = "OM1.17"
$ws.Range("C12").Value = "OM1.17"

$ws.Range("B13").Value = "OM1.18, OM1.42, OM1.46, OM7.3"

$ws.Range("D4").Value = "OM1.ProducerObservation ID"
$ws.Range("D8").Value = "OM1.ObservationChangeDateTime"
$ws.Range("D9").Value = "OM1.ObservationProducingDepartment"
$ws.Range("D11").Value = "OM1.ObservationPreferredLongName"
$ws.Range("D12").Value = "OM1.TelephoneNumber"
$ws.Range("D13").Value = "OM1.NatureOfObservation + OM1.KindOfQuantity + OM1.TargetAnatomicSiteOfTest + OM7. CategoryIdentifier"
$ws.Range("D6").Value = "MFE.RecordLevelEventCode"

$ws.Range("C13").Value = "OM1.18 + OM1.42 + OM1.46 + OM7.3"

# Leave the selection on the last edited cell, matching the author's
# final cursor position.
$ws.Range("B13").Select()
